# Update the remote DB IP used by the "database" sheet and make that
# sheet the active/selected tab (it was previously "simulation").

$wb = $excel.ActiveWorkbook

$dbSheet = $wb.Worksheets.Item("database")

# Change the configured IP value from "localhost" to the deployment IP.
$dbSheet.Range("A2").Value = "10.10.2.42"

# Make "database" the active tab instead of "simulation".
$dbSheet.Activate()
